# Apply the "S&P-500" benchmark block (columns N:S) to the existing
# "Activo" / "Estrategia" stats table and refresh the four data rows
# with the updated figures (eeuuAnalysis.xlsx work).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Header block (row 1): new merged "S&P-500" title in N1:S1,
#    cloned from the look of the existing H1:M1 ("Estrategia") block.
# ---------------------------------------------------------------
$ws.Range("H1:M1").Copy()
$ws.Range("N1:S1").PasteSpecial(-4122)
$ws.Range("N1:S1").Merge()
$ws.Range("N1").Value = "S&P-500"

# ---------------------------------------------------------------
# 2. Column headers (row 2): Mediana / Media geo. / Media arit. /
#    Desvio / Max / Min repeated for the new block.
# ---------------------------------------------------------------
$ws.Range("H2:M2").Copy()
$ws.Range("N2:S2").PasteSpecial(-4122)
$ws.Range("N2").Value = "Mediana"
$ws.Range("O2").Value = "Media geo."
$ws.Range("P2").Value = "Media arit."
$ws.Range("Q2").Value = "Desvio"
$ws.Range("R2").Value = "Max"
$ws.Range("S2").Value = "Min"

# ---------------------------------------------------------------
# 3. Data rows 4-7: refreshed "Activo"/"Estrategia" figures plus the
#    new "S&P-500" figures in columns N:S.
# ---------------------------------------------------------------

# periodo entero
$row4 = @(0.07000000000000001, 0.09, 0.11, 2.08, 11.98, -12.86, 0, 0.09, 0.11, 2.02, 12.86, -11.98, 0.08, 0.04, 0.05, 1.42, 9.380000000000001, -11.98)
$col = 2
foreach ($v in $row4) {
    $ws.Cells.Item(4, $col).Value = $v
    $col = $col + 1
}

# comprado
$row5 = @(0.08, 0.18, 0.19, 1.68, 10.47, -8.01, 0.08, 0.18, 0.19, 1.68, 10.47, -8.01, 0.15, 0.08, 0.08, 0.99, 3.41, -5.89)
$col = 2
foreach ($v in $row5) {
    $ws.Cells.Item(5, $col).Value = $v
    $col = $col + 1
}

# en efectivo
$row6 = @(0.37, 0.19, 0.21, 2.02, 7.56, -3.86, 0, 0, 0, 0, 0, 0, -0.08, 0.06, 0.06, 1.22, 2.76, -2.5)
$col = 2
foreach ($v in $row6) {
    $ws.Cells.Item(6, $col).Value = $v
    $col = $col + 1
}

# vendido
$row7 = @(0, -0.06, -0.02, 2.53, 11.98, -12.86, 0, -0.01, 0.02, 2.53, 12.86, -11.98, -0.05, -0.01, 0.01, 1.88, 9.380000000000001, -11.98)
$col = 2
foreach ($v in $row7) {
    $ws.Cells.Item(7, $col).Value = $v
    $col = $col + 1
}
